$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "4hr" note text (row 3, column D) to "4hr15min"
$ws.Range("D3").Value = "4hr15min"

# Fill in row 4 data: Time Cost -> 4, Time spent -> "0.5hr", Notes -> DONE text
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "0.5hr"
$ws.Range("E4").Value = "DONE Add <h2> with anchors; Style menu; JS autopopulate menu feature; Scroll between sections; NEW FORK: Tentative: inject aside with JS, inject anchors by scrubbing <h2> text nodes"

# Row height for row 4 to fit the wrapped note text
$ws.Rows.Item(4).RowHeight = 72.5

# Adjust column E width (stored width=39 once Excel re-derives it from this
# ColumnWidth) and select E4 as active cell, matching the recorded view state
$ws.Columns.Item(5).ColumnWidth = 38.14
$ws.Range("E4").Select()
